# Updated symbol list refresh (coinranking data pull):
#  - Several "Price" (column D) values refreshed with new quotes.
#  - Rows 10-18 (Coin / Link / Price / Volume) re-shuffled: a new "One" entry
#    pushed onto the list, cascading WazirX / LiechtensteinCryptoassetsExchange /
#    MandalaExchangeToken / BitrueCoin / MCDex / BitMartToken / BitForexToken /
#    CoinExToken down one row each, with fresh Price/Volume text for every row.
#  - Row 47's Volume(1h) label lost its "Bestin24h" suffix.
#
# Column D historically stores the price as literal TEXT (e.g. "245.88"), not
# a number, so plain numeric-looking assignments are forced back to text with
# a leading apostrophe (Excel's "store as text" quote-prefix) and the cell's
# style is then reset to Normal so no stray number-format/quote-prefix style
# lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $value) {
    $cell = $ws.Range($a1)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# --- Simple price-only (column D) refreshes ---
Set-TextValue "D2" "245.98"
Set-TextValue "D4" "5.296"
Set-TextValue "D5" "0.05875"
Set-TextValue "D6" "3.383"
Set-TextValue "D7" "6.377"
Set-TextValue "D8" "0.8170"
Set-TextValue "D9" "0.9614"

# --- Rows 10-18: coin/link/price/volume block re-shuffled ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01120"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1418"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03599"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D13" "0.07319"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03041"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "4.445"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09386"
$ws.Range("E16").Value = "15BitMartTokenBMX"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001599"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04820"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- More standalone price (column D) refreshes ---
Set-TextValue "D19" "0.006187"
Set-TextValue "D21" "0.0009882"
Set-TextValue "D22" "0.00009703"
Set-TextValue "D23" "3.682"
Set-TextValue "D25" "0.3262"
Set-TextValue "D26" "0.1277"
Set-TextValue "D40" "0.03862"
Set-TextValue "D41" "0.006602"
Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.005887"
Set-TextValue "D45" "0.00005663"

Set-TextValue "D47" "0.7753"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.08535"
